$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Helper: write a text value into a cell while preserving the cells original
# (unstyled) appearance. Plain assignment to .Value lets Excel auto-detect the
# string as a number, which would corrupt values such as "26.900.39" or drop
# significant trailing zeros (e.g. "0.3670" -> 0.367). Temporarily forcing the
# cell to Text format makes Excel keep the literal string, and ClearFormats()
# afterwards removes the now-unneeded style index again so the cell keeps the
# same (default) style it had before the edit.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextValue ($ws.Cells.Item(2, 4)) "26.900.39"
Set-TextValue ($ws.Cells.Item(2, 5)) "  -2.23%  "
Set-TextValue ($ws.Cells.Item(3, 4)) "1.835.10"
Set-TextValue ($ws.Cells.Item(3, 5)) "  -1.66%  "
Set-TextValue ($ws.Cells.Item(4, 4)) "1.006"
Set-TextValue ($ws.Cells.Item(4, 5)) "  +0.17%  "
Set-TextValue ($ws.Cells.Item(5, 4)) "310.55"
Set-TextValue ($ws.Cells.Item(5, 5)) "  -1.70%  "
Set-TextValue ($ws.Cells.Item(6, 5)) "  +0.17%  "
Set-TextValue ($ws.Cells.Item(7, 4)) "0.4617"
Set-TextValue ($ws.Cells.Item(7, 5)) "  -1.35%  "
Set-TextValue ($ws.Cells.Item(8, 4)) "0.3670"
Set-TextValue ($ws.Cells.Item(8, 5)) "  -1.61%  "
Set-TextValue ($ws.Cells.Item(9, 4)) "0.07178"
Set-TextValue ($ws.Cells.Item(9, 5)) "  -2.80%  "
Set-TextValue ($ws.Cells.Item(10, 4)) "0.8808"
Set-TextValue ($ws.Cells.Item(10, 5)) "  -0.96%  "
Set-TextValue ($ws.Cells.Item(11, 4)) "0.07856"
Set-TextValue ($ws.Cells.Item(11, 5)) "  -0.92%  "
Set-TextValue ($ws.Cells.Item(12, 5)) "  -1.97%  "
Set-TextValue ($ws.Cells.Item(13, 4)) "1.860.20"
Set-TextValue ($ws.Cells.Item(13, 5)) "  -0.29%  "
Set-TextValue ($ws.Cells.Item(14, 4)) "5.356"
Set-TextValue ($ws.Cells.Item(14, 5)) "  -1.21%  "
Set-TextValue ($ws.Cells.Item(15, 4)) "6.408"
Set-TextValue ($ws.Cells.Item(15, 5)) "  -2.95%  "
Set-TextValue ($ws.Cells.Item(16, 4)) "88.39"
Set-TextValue ($ws.Cells.Item(16, 5)) "  -4.72%  "
Set-TextValue ($ws.Cells.Item(17, 5)) "  +0.13%  "
Set-TextValue ($ws.Cells.Item(18, 4)) "0.000008746"
Set-TextValue ($ws.Cells.Item(18, 5)) "  -2.15%  "
Set-TextValue ($ws.Cells.Item(19, 5)) "  +0.13%  "
Set-TextValue ($ws.Cells.Item(20, 4)) "26.934.47"
Set-TextValue ($ws.Cells.Item(20, 5)) "  -2.24%  "
Set-TextValue ($ws.Cells.Item(21, 4)) "14.50"
Set-TextValue ($ws.Cells.Item(21, 5)) "  -2.68%  "
Set-TextValue ($ws.Cells.Item(22, 4)) "5.020"
Set-TextValue ($ws.Cells.Item(22, 5)) "  -2.66%  "
Set-TextValue ($ws.Cells.Item(23, 5)) "  -1.24%  "
Set-TextValue ($ws.Cells.Item(24, 4)) "1.981"
Set-TextValue ($ws.Cells.Item(24, 5)) "  +4.97%  "
Set-TextValue ($ws.Cells.Item(25, 4)) "151.07"
Set-TextValue ($ws.Cells.Item(25, 5)) "  -1.57%  "
Set-TextValue ($ws.Cells.Item(26, 4)) "18.28"
Set-TextValue ($ws.Cells.Item(26, 5)) "  -1.44%  "
Set-TextValue ($ws.Cells.Item(27, 5)) "  -4.97%  "
Set-TextValue ($ws.Cells.Item(28, 4)) "113.97"
Set-TextValue ($ws.Cells.Item(28, 5)) "  -2.54%  "
Set-TextValue ($ws.Cells.Item(29, 4)) "4.958"
Set-TextValue ($ws.Cells.Item(29, 5)) "  -4.31%  "
Set-TextValue ($ws.Cells.Item(30, 4)) "0.08849"
Set-TextValue ($ws.Cells.Item(30, 5)) "  -0.66%  "
Set-TextValue ($ws.Cells.Item(31, 4)) "3.138"
Set-TextValue ($ws.Cells.Item(31, 5)) "  +3.71%  "
Set-TextValue ($ws.Cells.Item(32, 4)) "0.7603"
Set-TextValue ($ws.Cells.Item(32, 5)) "  +0.32%  "
Set-TextValue ($ws.Cells.Item(33, 4)) "4.477"
Set-TextValue ($ws.Cells.Item(33, 5)) "  -0.32%  "
Set-TextValue ($ws.Cells.Item(34, 4)) "1.132"
Set-TextValue ($ws.Cells.Item(34, 5)) "  -3.25%  "
Set-TextValue ($ws.Cells.Item(35, 4)) "2.644"
Set-TextValue ($ws.Cells.Item(35, 5)) "  +0.69%  "
Set-TextValue ($ws.Cells.Item(36, 4)) "1.090"
Set-TextValue ($ws.Cells.Item(36, 5)) "  +0.72%  "
Set-TextValue ($ws.Cells.Item(37, 4)) "0.01938"
Set-TextValue ($ws.Cells.Item(37, 5)) "  -1.55%  "
Set-TextValue ($ws.Cells.Item(38, 4)) "2.931"
Set-TextValue ($ws.Cells.Item(38, 5)) "  -1.98%  "
Set-TextValue ($ws.Cells.Item(39, 4)) "0.05146"
Set-TextValue ($ws.Cells.Item(39, 5)) "  -2.55%  "
Set-TextValue ($ws.Cells.Item(40, 5)) "  -3.43%  "
Set-TextValue ($ws.Cells.Item(41, 4)) "0.4990"
Set-TextValue ($ws.Cells.Item(41, 5)) "  -4.34%  "
Set-TextValue ($ws.Cells.Item(42, 4)) "0.1601"
Set-TextValue ($ws.Cells.Item(42, 5)) "  -2.80%  "
Set-TextValue ($ws.Cells.Item(43, 4)) "8.321"
Set-TextValue ($ws.Cells.Item(43, 5)) "  -0.82%  "
Set-TextValue ($ws.Cells.Item(44, 4)) "0.4705"
Set-TextValue ($ws.Cells.Item(44, 5)) "  -3.59%  "
Set-TextValue ($ws.Cells.Item(45, 2)) "PaxDollar"
Set-TextValue ($ws.Cells.Item(45, 3)) "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue ($ws.Cells.Item(45, 4)) "1.006"
Set-TextValue ($ws.Cells.Item(45, 5)) "  +0.20%  "
Set-TextValue ($ws.Cells.Item(46, 2)) "EnergySwap"
Set-TextValue ($ws.Cells.Item(46, 3)) "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue ($ws.Cells.Item(46, 4)) "10.13"
Set-TextValue ($ws.Cells.Item(46, 5)) "  -1.88%  "
Set-TextValue ($ws.Cells.Item(47, 4)) "102.68"
Set-TextValue ($ws.Cells.Item(47, 5)) "  -1.26%  "
Set-TextValue ($ws.Cells.Item(48, 4)) "1.615"
Set-TextValue ($ws.Cells.Item(48, 5)) "  -2.50%  "
Set-TextValue ($ws.Cells.Item(49, 5)) "  -2.83%  "
Set-TextValue ($ws.Cells.Item(50, 4)) "64.94"
Set-TextValue ($ws.Cells.Item(50, 5)) "  -1.48%  "
Set-TextValue ($ws.Cells.Item(51, 4)) "36.46"
Set-TextValue ($ws.Cells.Item(51, 5)) "  -2.01%  "
